$d = $word.ActiveDocument

# Locate the existing trailing empty paragraph that carries <w:ind w:left="360"/>.
# It is the second-to-last paragraph in the document body (the very last
# paragraph is the final, completely blank one right before the sectPr).
$anchorIndex = $d.Paragraphs.Count - 1
$anchor = $d.Paragraphs.Item($anchorIndex)
$insertionPoint = $d.Range($anchor.Range.Start, $anchor.Range.Start)

$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Build the new content as raw WordprocessingML so formatting (bold runs,
# ListParagraph/numPr bullet items, multiple runs in one paragraph) comes
# out exactly as desired. A trailing empty <w:p> is appended so that the
# last real paragraph gets its own paragraph mark instead of being fused
# into the pre-existing anchor paragraph; that spare paragraph mark is
# removed again right afterwards.
$xml = @"
<w:p $w><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>
<w:p $w><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>5/ OS:</w:t></w:r></w:p>
<w:p $w><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Window</w:t></w:r></w:p>
<w:p $w><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Linux</w:t></w:r></w:p>
<w:p $w><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Mac OS</w:t></w:r></w:p>
<w:p $w><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>It will provide service to the computer</w:t></w:r><w:r><w:t xml:space="preserve"> with the resource.</w:t></w:r></w:p>
<w:p $w></w:p>
"@

[void]$insertionPoint.InsertXML($xml)

# Remove the spare paragraph mark introduced above so the original anchor
# paragraph (still holding <w:ind w:left="360"/>) is left completely intact.
$spareIndex = $anchorIndex + 6
$spare = $d.Paragraphs.Item($spareIndex)
$spareMark = $d.Range($spare.Range.End - 1, $spare.Range.End)
$spareMark.Delete()
